$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("Förändrad") bumps from 46059 to 46060 for rows 2-9 ---
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value2 = 46060
}

# --- Rows 4-9 content rotates: rows 7,8,9 move up to 4,5,6 and rows 4,5,6 move down to 7,8,9 ---
# Capture the current (pre-edit) values for rows 4-9 in columns A, B, F, G
$data = @{}
foreach ($r in 4..9) {
    $data[$r] = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
    }
}

# Mapping: new row <- old row
$mapping = @{ 4 = 7; 5 = 8; 6 = 9; 7 = 4; 8 = 5; 9 = 6 }

foreach ($newRow in 4..9) {
    $oldRow = $mapping[$newRow]
    $src = $data[$oldRow]

    $ws.Cells.Item($newRow, 1).Value = $src.A

    $ws.Cells.Item($newRow, 2).Value2 = $src.B

    if ($src.F) {
        $ws.Cells.Item($newRow, 6).Value = $src.F
    } else {
        $ws.Cells.Item($newRow, 6).ClearContents()
    }

    $ws.Cells.Item($newRow, 7).Value = $src.G
}
